# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 593.8333
$ws.Cells.Item(2, 9).Value = 90.25
$ws.Cells.Item(2, 10).Value = 1601
$ws.Cells.Item(2, 11).Value = 90.25
$ws.Cells.Item(2, 12).Value = 1601
$ws.Cells.Item(2, 13).Value = 22.75
$ws.Cells.Item(2, 14).Value = -1827

$ws.Cells.Item(18, 8).Value = 10000
$ws.Cells.Item(18, 9).Value = 10000
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 10000
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -9716
$ws.Cells.Item(18, 14).ClearContents()

$ws.Cells.Item(33, 8).Value = 196.5
$ws.Cells.Item(33, 9).Value = 230
$ws.Cells.Item(33, 11).Value = 230
$ws.Cells.Item(33, 13).Value = -1

$ws.Cells.Item(40, 8).Value = 2487.75
$ws.Cells.Item(40, 9).Value = 2380.2
$ws.Cells.Item(40, 10).Value = 2667
$ws.Cells.Item(40, 11).Value = 2380.2
$ws.Cells.Item(40, 12).Value = 2667
$ws.Cells.Item(40, 13).Value = -2205.2
$ws.Cells.Item(40, 14).Value = -3017

$ws.Cells.Item(101, 8).Value = 664
$ws.Cells.Item(101, 9).Value = 283.75
$ws.Cells.Item(101, 10).Value = 2185
$ws.Cells.Item(101, 11).Value = 851.25
$ws.Cells.Item(101, 12).Value = 6555
$ws.Cells.Item(101, 13).Value = 770.75
$ws.Cells.Item(101, 14).Value = -9799

$ws.Cells.Item(112, 8).Value = 4833.3335
$ws.Cells.Item(112, 10).Value = 4833.3335
$ws.Cells.Item(112, 12).Value = 14500.0005
$ws.Cells.Item(112, 14).Value = -16716.0005

$ws.Cells.Item(118, 8).Value = 757.3
$ws.Cells.Item(118, 9).Value = 757.3
$ws.Cells.Item(118, 11).Value = 2271.9
$ws.Cells.Item(118, 13).Value = -614.8999999999996

$ws.Cells.Item(138, 8).Value = 2096.5
$ws.Cells.Item(138, 9).Value = 2096.5
$ws.Cells.Item(138, 11).Value = 6289.5
$ws.Cells.Item(138, 13).Value = -1149.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 14).ClearContents()

$ws.Cells.Item(61, 8).Value = 3999.889
$ws.Cells.Item(61, 9).Value = 3499.875
$ws.Cells.Item(61, 11).Value = 3499.875
$ws.Cells.Item(61, 13).Value = -3287.875

$ws.Cells.Item(110, 8).Value = 5325.9165
$ws.Cells.Item(110, 9).Value = 5630.1177
$ws.Cells.Item(110, 10).Value = 4587.143
$ws.Cells.Item(110, 11).Value = 5630.1177
$ws.Cells.Item(110, 12).Value = 4587.143
$ws.Cells.Item(110, 13).Value = -3585.1177
$ws.Cells.Item(110, 14).Value = -8677.143

$ws.Cells.Item(136, 8).Value = 3999.889
$ws.Cells.Item(136, 9).Value = 3499.875
$ws.Cells.Item(136, 11).Value = 10499.625
$ws.Cells.Item(136, 13).Value = -7949.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(6, 8).Value = 10749.417
$ws.Cells.Item(6, 10).Value = 10749.417
$ws.Cells.Item(6, 12).Value = 10749.417
$ws.Cells.Item(6, 14).Value = -10975.417

$ws.Cells.Item(75, 8).Value = 4995
$ws.Cells.Item(75, 9).Value = 4995
$ws.Cells.Item(75, 11).Value = 4995
$ws.Cells.Item(75, 13).Value = -4059

$ws.Cells.Item(78, 8).Value = 4995
$ws.Cells.Item(78, 9).Value = 4995
$ws.Cells.Item(78, 11).Value = 14985
$ws.Cells.Item(78, 13).Value = -10305

$ws.Cells.Item(86, 8).Value = 7289.4
$ws.Cells.Item(86, 9).Value = 3413.7144
$ws.Cells.Item(86, 11).Value = 3413.7144
$ws.Cells.Item(86, 13).Value = -2290.7144

$ws.Cells.Item(89, 8).Value = 7289.4
$ws.Cells.Item(89, 9).Value = 3413.7144
$ws.Cells.Item(89, 11).Value = 17068.572
$ws.Cells.Item(89, 13).Value = -11452.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2253.1667
$ws.Cells.Item(31, 9).Value = 2130
$ws.Cells.Item(31, 10).Value = 2499.5
$ws.Cells.Item(31, 11).Value = 2130
$ws.Cells.Item(31, 12).Value = 2499.5
$ws.Cells.Item(31, 13).Value = -1835
$ws.Cells.Item(31, 14).Value = -3089.5

$ws.Cells.Item(34, 8).Value = 2253.1667
$ws.Cells.Item(34, 9).Value = 2130
$ws.Cells.Item(34, 10).Value = 2499.5
$ws.Cells.Item(34, 11).Value = 2130
$ws.Cells.Item(34, 12).Value = 2499.5
$ws.Cells.Item(34, 13).Value = -1928
$ws.Cells.Item(34, 14).Value = -2903.5

$ws.Cells.Item(60, 8).Value = 56746.668
$ws.Cells.Item(60, 9).Value = 75000
$ws.Cells.Item(60, 10).Value = 51531.43
$ws.Cells.Item(60, 11).Value = 75000
$ws.Cells.Item(60, 12).Value = 51531.43
$ws.Cells.Item(60, 13).Value = -74489
$ws.Cells.Item(60, 14).Value = -52553.43

$ws.Cells.Item(94, 8).Value = 1569.4
$ws.Cells.Item(94, 9).Value = 1636.75
$ws.Cells.Item(94, 10).Value = 1300
$ws.Cells.Item(94, 11).Value = 1636.75
$ws.Cells.Item(94, 12).Value = 1300
$ws.Cells.Item(94, 13).Value = -1185.75
$ws.Cells.Item(94, 14).Value = -2202

$ws.Cells.Item(121, 8).Value = 59576
$ws.Cells.Item(121, 10).Value = 59576
$ws.Cells.Item(121, 12).Value = 59576
$ws.Cells.Item(121, 14).Value = -62196

$ws.Cells.Item(131, 8).Value = 60000
$ws.Cells.Item(131, 10).Value = 60000
$ws.Cells.Item(131, 12).Value = 60000
$ws.Cells.Item(131, 14).Value = -70080

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 929
$ws.Cells.Item(86, 10).Value = 1859.75
$ws.Cells.Item(86, 12).Value = 5579.25
$ws.Cells.Item(86, 14).Value = -7951.25

$ws.Cells.Item(89, 8).Value = 929
$ws.Cells.Item(89, 10).Value = 1859.75
$ws.Cells.Item(89, 12).Value = 16737.75
$ws.Cells.Item(89, 14).Value = -28593.75

$ws.Cells.Item(100, 8).Value = 3999.9167
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 3999.9167
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 11999.7501
$ws.Cells.Item(100, 13).ClearContents()
$ws.Cells.Item(100, 14).Value = -13621.7501

$ws.Cells.Item(106, 8).Value = 12760.667
$ws.Cells.Item(106, 10).Value = 13263.077
$ws.Cells.Item(106, 12).Value = 39789.231
$ws.Cells.Item(106, 14).Value = -41681.231

$ws.Cells.Item(131, 8).Value = 1976.2
$ws.Cells.Item(131, 10).Value = 1970.25
$ws.Cells.Item(131, 12).Value = 5910.75
$ws.Cells.Item(131, 14).Value = -15990.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 52059
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 13).ClearContents()

$ws.Cells.Item(81, 8).Value = 52059
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 13).ClearContents()

$ws.Cells.Item(84, 8).Value = 52059
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 13).ClearContents()

$ws.Cells.Item(102, 8).Value = 3529.4
$ws.Cells.Item(102, 9).Value = 3588.2222
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 3588.2222
$ws.Cells.Item(102, 12).Value = 3000
$ws.Cells.Item(102, 13).Value = -1966.2222
$ws.Cells.Item(102, 14).Value = -6244

$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 2749.5
$ws.Cells.Item(132, 9).Value = 2749.5
$ws.Cells.Item(132, 11).Value = 8248.5
$ws.Cells.Item(132, 13).Value = -5718.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3878.5386
$ws.Cells.Item(7, 9).Value = 3606
$ws.Cells.Item(7, 10).Value = 4314.6
$ws.Cells.Item(7, 11).Value = 3606
$ws.Cells.Item(7, 12).Value = 4314.6
$ws.Cells.Item(7, 13).Value = -3494
$ws.Cells.Item(7, 14).Value = -4538.6

$ws.Cells.Item(46, 8).Value = 1677.4445
$ws.Cells.Item(46, 9).Value = 974
$ws.Cells.Item(46, 11).Value = 974
$ws.Cells.Item(46, 13).Value = -786

$ws.Cells.Item(126, 8).Value = 3878.5386
$ws.Cells.Item(126, 9).Value = 3606
$ws.Cells.Item(126, 10).Value = 4314.6
$ws.Cells.Item(126, 11).Value = 10818
$ws.Cells.Item(126, 12).Value = 12943.8
$ws.Cells.Item(126, 13).Value = -8348
$ws.Cells.Item(126, 14).Value = -17883.8

$ws.Cells.Item(132, 8).Value = 1985.2727
$ws.Cells.Item(132, 9).Value = 1964.4166
$ws.Cells.Item(132, 10).Value = 2040.8889
$ws.Cells.Item(132, 11).Value = 5893.2498
$ws.Cells.Item(132, 12).Value = 6122.6667
$ws.Cells.Item(132, 13).Value = -3363.2498
$ws.Cells.Item(132, 14).Value = -11182.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3385.7856
$ws.Cells.Item(96, 9).Value = 2495.8572
$ws.Cells.Item(96, 10).Value = 4275.7144
$ws.Cells.Item(96, 11).Value = 2495.8572
$ws.Cells.Item(96, 12).Value = 4275.7144
$ws.Cells.Item(96, 13).Value = -1122.8572
$ws.Cells.Item(96, 14).Value = -7021.7144

$ws.Cells.Item(136, 8).Value = 6896.12
$ws.Cells.Item(136, 9).Value = 6896.12
$ws.Cells.Item(136, 11).Value = 20688.36
$ws.Cells.Item(136, 13).Value = -18138.36
